# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) previously held a different statistic ("Strike#"); it is
# regenerated here with the new per-row K values so the sheet reflects the
# regenerated save_data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 0
    5  = 1
    6  = 0
    7  = 3
    8  = 2
    9  = 3
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 2
    16 = 1
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 0
    23 = 0
    24 = 2
    25 = 0
    26 = 0
    27 = 0
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 0
    33 = 3
    34 = 1
    35 = 1
    36 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
